$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.96249999999999
$ws.Range("A6").Value = -22.92410000000001
$ws.Range("A7").Value = -21.9234
$ws.Range("B7").Value = 5.023699999999999
$ws.Range("A8").Value = -22.33810000000001
$ws.Range("B11").Value = 5.975799999999998
$ws.Range("B12").Value = 5.4396
$ws.Range("C12").Value = -10.97599999999999
$ws.Range("C13").Value = -13.31559999999999
$ws.Range("C14").Value = -14.51939999999999
$ws.Range("B15").Value = 5.058999999999997
$ws.Range("A16").Value = -21.57529999999999
$ws.Range("C16").Value = -11.68999999999999
$ws.Range("C19").Value = -12.19030000000001
$ws.Range("A20").Value = -22.4795
$ws.Range("B20").Value = 4.419299999999996
$ws.Range("C20").Value = -14.28449999999999
$ws.Range("A21").Value = -22.48799999999999
$ws.Range("B21").Value = 5.089399999999999
$ws.Range("B22").Value = 10.0382
$ws.Range("C22").Value = -12.392
$ws.Range("B23").Value = 9.6774
$ws.Range("A28").Value = -22.02209999999999
$ws.Range("A29").Value = -21.69030000000001
$ws.Range("B29").Value = 4.923400000000004
$ws.Range("A30").Value = -21.73500000000002
$ws.Range("A32").Value = -21.42550000000001
$ws.Range("B34").Value = 9.152600000000005
$ws.Range("C36").Value = -12.59740000000001
$ws.Range("A40").Value = -19.65749999999998
$ws.Range("B42").Value = 9.847499999999995
$ws.Range("B43").Value = 4.994199999999995
$ws.Range("C43").Value = -13.4171
$ws.Range("B44").Value = 5.045100000000001
$ws.Range("B45").Value = 5.043500000000001
$ws.Range("A46").Value = -22.284
$ws.Range("B46").Value = 5.307099999999994
$ws.Range("C46").Value = -13.664
$ws.Range("B50").Value = 4.834199999999996
$ws.Range("C50").Value = -13.96989999999999
$ws.Range("A51").Value = -22.36489999999999
$ws.Range("B51").Value = 4.8528
$ws.Range("A52").Value = -22.1006
$ws.Range("A57").Value = -22.78230000000001
$ws.Range("B57").Value = 5.003399999999996
$ws.Range("A59").Value = -22.24810000000001
$ws.Range("A62").Value = -22.32120000000001
$ws.Range("B65").Value = 5.282500000000002
$ws.Range("A66").Value = -21.5317
$ws.Range("B66").Value = 4.807599999999997
$ws.Range("B67").Value = 5.233799999999999
$ws.Range("A73").Value = -19.67839999999998
$ws.Range("A74").Value = -22.10379999999999
$ws.Range("C76").Value = -12.0539
$ws.Range("A77").Value = -20.37699999999997
$ws.Range("B79").Value = 9.715000000000007
$ws.Range("B84").Value = 5.712900000000001
$ws.Range("B87").Value = 5.090399999999999
$ws.Range("A92").Value = -21.36510000000001
$ws.Range("B92").Value = 5.915599999999996
$ws.Range("C95").Value = -11.58870000000001
$ws.Range("B97").Value = 6.238499999999996
$ws.Range("C97").Value = -11.16560000000001
$ws.Range("C99").Value = -12.2339
$ws.Range("A100").Value = -22.3255
